$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$questionsText = @'
questions = [
    {
        "title": "Your organization plans to expand its Azure-based services globally, requiring a resilient virtual network design to ensure uninterrupted services across multiple Azure regions. You are tasked with designing a virtual network architecture that facilitates low-latency, secure, and reliable interconnectivity.Which combination of Azure services should you use?",
        "ques_type": 2,
        "options": [
            "Azure Virtual WAN &amp Azure ExpressRoute",
            "Azure Virtual Network &amp Azure VPN Gateway",
            "Azure Load Balancer &amp Azure Application Gateway",
            "Azure Traffic Manager &amp Azure Bastion"
        ],
        "score": "Azure Virtual WAN &amp Azure ExpressRoute"
    },
    {
        "title": "Your company's Azure infrastructure is experiencing an increase in unauthorized access attempts. You need to implement a solution that not only blocks these attempts but also provides advanced threat protection and integrated security management.Which combination of Azure services should you use?",
        "ques_type": 2,
        "options": [
            "Azure Firewall &amp Azure Security Center",
            "Azure Application Gateway with Web Application Firewall (WAF) &amp Azure Sentinel",
            "Azure Bastion &amp Microsoft Entra ID",
            "Azure VPN Gateway &amp Azure Monitor"
        ],
        "score": "Azure Firewall &amp Azure Security Center"
    },
    {
        "title": "You are responsible for maintaining the network health of your Azure environment. You need to set up a solution that allows real-time monitoring, advanced analytics, and the ability to respond to network performance issues.Which Azure service should you primarily use?",
        "ques_type": 2,
        "options": [
            "Azure Network Watcher",
            "Azure Monitor",
            "Azure Application Insights",
            "Azure Service Health"
        ],
        "score": "Azure Network Watcher"
    },
    {
        "title": "Your organization operates a hybrid cloud environment with Azure and on-premises data centers. You need a solution that seamlessly connects on-premises networks to Azure, ensuring secure and reliable communication.Which Azure service should you use?",
        "ques_type": 2,
        "options": [
            "Azure ExpressRoute",
            "Azure Virtual Network",
            "Azure VPN Gateway",
            "Azure Virtual WAN"
        ],
        "score": "Azure ExpressRoute"
    }
]
'@
$ws.Range("A2").ClearContents()
$ws.Range("A1").Value = $questionsText
$ws.Range("A1").Style = "Normal"
